$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 3569
$ws.Range("I3").Value = 3693
$ws.Range("I4").Value = 869
$ws.Range("I5").Value = 343
$ws.Range("I6").Value = 4137
$ws.Range("I7").Value = 12611

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I2").Value = 34
$ws.Range("I6").Value = 53
$ws.Range("I7").Value = 143

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I6").Value = 54
$ws.Range("I7").Value = 142

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I2").Value = 18
$ws.Range("I7").Value = 62

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I3").Value = 121
$ws.Range("I6").Value = 115
$ws.Range("I7").Value = 404

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I4").Value = 18
$ws.Range("I7").Value = 237

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 117
$ws.Range("I6").Value = 164
$ws.Range("I7").Value = 486

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I2").Value = 38
$ws.Range("I7").Value = 118

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I2").Value = 40
$ws.Range("I7").Value = 110

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I6").Value = 88
$ws.Range("I7").Value = 410
$ws.Range("I8").Value = 767
$ws.Range("I11").Value = 201
$ws.Range("I12").Value = 26
$ws.Range("I13").Value = 21
$ws.Range("I14").Value = 62
$ws.Range("I16").Value = 34
$ws.Range("I19").Value = 332
$ws.Range("I20").Value = 313
$ws.Range("I23").Value = 118
$ws.Range("I25").Value = 60
$ws.Range("I26").Value = 18
$ws.Range("I27").Value = 113
$ws.Range("I29").Value = 826
$ws.Range("I31").Value = 118
$ws.Range("I33").Value = 568
$ws.Range("I34").Value = 58
$ws.Range("I36").Value = 174
$ws.Range("I37").Value = 404
$ws.Range("I42").Value = 437
$ws.Range("I47").Value = 88
$ws.Range("I49").Value = 104
$ws.Range("I54").Value = 282
$ws.Range("I59").Value = 25
$ws.Range("I63").Value = 48
$ws.Range("I67").Value = 486
$ws.Range("I73").Value = 110
$ws.Range("I74").Value = 28
$ws.Range("I76").Value = 193
$ws.Range("I78").Value = 181
$ws.Range("I79").Value = 329
$ws.Range("I84").Value = 110
$ws.Range("I85").Value = 579
$ws.Range("I86").Value = 77
$ws.Range("I89").Value = 143
$ws.Range("I90").Value = 160
$ws.Range("I91").Value = 152
$ws.Range("I93").Value = 70
$ws.Range("I94").Value = 113
$ws.Range("I95").Value = 199
$ws.Range("I96").Value = 142
$ws.Range("I98").Value = 82
$ws.Range("I99").Value = 237
$ws.Range("I101").Value = 12611

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I2").Value = 71
$ws.Range("I7").Value = 199

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 135
$ws.Range("I6").Value = 182
$ws.Range("I7").Value = 568

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("I6").Value = 66
$ws.Range("I7").Value = 104

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I6").Value = 144
$ws.Range("I7").Value = 282

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 249
$ws.Range("I3").Value = 283
$ws.Range("I6").Value = 225
$ws.Range("I7").Value = 826

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I6").Value = 92
$ws.Range("I7").Value = 332

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I4").Value = 22
$ws.Range("I7").Value = 193

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I6").Value = 145
$ws.Range("I7").Value = 579

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I2").Value = 39
$ws.Range("I7").Value = 88

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I3").Value = 150
$ws.Range("I7").Value = 437

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range("I4").Value = 8
$ws.Range("I6").Value = 21

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I6").Value = 74
$ws.Range("I7").Value = 181

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I2").Value = 32
$ws.Range("I7").Value = 118

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I6").Value = 47
$ws.Range("I7").Value = 152

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I3").Value = 109
$ws.Range("I7").Value = 329

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I3").Value = 94
$ws.Range("I4").Value = 21
$ws.Range("I7").Value = 313

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I2").Value = 56
$ws.Range("I3").Value = 53
$ws.Range("I6").Value = 52
$ws.Range("I7").Value = 174

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("I2").Value = 21
$ws.Range("I7").Value = 70

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("I3").Value = 20
$ws.Range("I7").Value = 58

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I2").Value = 20
$ws.Range("I6").Value = 63
$ws.Range("I7").Value = 113

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("I2").Value = 18
$ws.Range("I3").Value = 16
$ws.Range("I7").Value = 60

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I3").Value = 28
$ws.Range("I6").Value = 34
$ws.Range("I7").Value = 88

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I3").Value = 8
$ws.Range("I6").Value = 50
$ws.Range("I7").Value = 82

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("I3").Value = 4
$ws.Range("I7").Value = 18

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I2").Value = 89
$ws.Range("I6").Value = 47
$ws.Range("I7").Value = 201

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("I3").Value = 30
$ws.Range("I4").Value = 12
$ws.Range("I7").Value = 110

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("I2").Value = 11
$ws.Range("I7").Value = 25

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 241
$ws.Range("I5").Value = 25
$ws.Range("I6").Value = 244
$ws.Range("I7").Value = 767

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("I2").Value = 29
$ws.Range("I6").Value = 46
$ws.Range("I7").Value = 113

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("I2").Value = 14
$ws.Range("I4").Value = 40
$ws.Range("I7").Value = 77

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I6").Value = 55
$ws.Range("I7").Value = 160

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("I2").Value = 7
$ws.Range("I3").Value = 10

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 140
$ws.Range("I3").Value = 122
$ws.Range("I7").Value = 410

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("I2").Value = 5
$ws.Range("I3").Value = 4
$ws.Range("I7").Value = 26

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("I2").Value = 8
$ws.Range("I7").Value = 34

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 28

$wb.Save()
Write-Host "Applied 160 cell updates across 46 sheets."